$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-04 Friday", "2025-07-05 Saturday"),
    @("50×39=", "21×36="),
    @("62×67=", "72×89="),
    @("74×56=", "82×88="),
    @("90×82=", "36×73="),
    @("33×27=", "54×63="),
    @("93×19=", "59×24="),
    @("53×82=", "17×18="),
    @("81×60=", "91×48="),
    @("52×32=", "32×39="),
    @("41×26=", "78×65="),
    @("30×86=", "95×47="),
    @("58×87=", "36×18="),
    @("31×27=", "26×32="),
    @("17×75=", "78×54="),
    @("73×46=", "58×26="),
    @("47×55=", "15×62="),
    @("51×54=", "79×40="),
    @("81×23=", "53×78="),
    @("95×31=", "98×29="),
    @("29×98=", "20×95="),
    @("77×37=", "74×50="),
    @("36×96=", "44×31="),
    @("41×21=", "26×53="),
    @("70×32=", "54×48="),
    @("54×92=", "59×99=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
